# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to Belias_Profits Leve tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 12307.444
$ws.Range("I21").Value = 6491.75
$ws.Range("J21").Value = 16960
$ws.Range("K21").Value = 6491.75
$ws.Range("L21").Value = 16960
$ws.Range("M21").Value = -6023.75
$ws.Range("N21").Value = -17896
# Row 23
$ws.Range("H23").Value = 12307.444
$ws.Range("I23").Value = 6491.75
$ws.Range("J23").Value = 16960
$ws.Range("K23").Value = 6491.75
$ws.Range("L23").Value = 16960
$ws.Range("M23").Value = -6257.75
$ws.Range("N23").Value = -17428
# Row 33
$ws.Range("H33").Value = 516.4286
$ws.Range("I33").Value = 224.30435
$ws.Range("K33").Value = 224.30435
$ws.Range("M33").Value = 4.695650000000001
# Row 63
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
# Row 66
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
# Row 113
$ws.Range("H113").Value = 5830.9473
$ws.Range("I113").Value = 4288.8887
$ws.Range("J113").Value = 7218.8
$ws.Range("K113").Value = 4288.8887
$ws.Range("L113").Value = 7218.8
$ws.Range("M113").Value = -1034.8887
$ws.Range("N113").Value = -13726.8
# Row 125
$ws.Range("H125").Value = 2597.318
$ws.Range("I125").Value = 2016.5
$ws.Range("J125").Value = 2929.2144
$ws.Range("K125").Value = 18148.5
$ws.Range("L125").Value = 26362.9296
$ws.Range("M125").Value = -15688.5
$ws.Range("N125").Value = -31282.9296
# Row 127
$ws.Range("H127").Value = 40000604
$ws.Range("I127").Value = 55555930
$ws.Range("J127").Value = 1200
$ws.Range("K127").Value = 166667790
$ws.Range("L127").Value = 3600
$ws.Range("M127").Value = -166662830
$ws.Range("N127").Value = -13520
# Row 138
$ws.Range("H138").Value = 2020.875
$ws.Range("I138").Value = 1357.7391
$ws.Range("J138").Value = 2918.0588
$ws.Range("K138").Value = 4073.2173
$ws.Range("L138").Value = 8754.1764
$ws.Range("M138").Value = 1066.7827
$ws.Range("N138").Value = -19034.1764

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6913.6665
$ws.Range("I61").Value = 1143.05
$ws.Range("J61").Value = 23401.143
$ws.Range("K61").Value = 1143.05
$ws.Range("L61").Value = 23401.143
$ws.Range("M61").Value = -931.05
$ws.Range("N61").Value = -23825.143
# Row 102
$ws.Range("H102").Value = 2345.5386
$ws.Range("I102").Value = 2763.75
$ws.Range("J102").Value = 1676.4
$ws.Range("K102").Value = 2763.75
$ws.Range("L102").Value = 1676.4
$ws.Range("M102").Value = -1141.75
$ws.Range("N102").Value = -4920.4
# Row 136
$ws.Range("H136").Value = 6913.6665
$ws.Range("I136").Value = 1143.05
$ws.Range("J136").Value = 23401.143
$ws.Range("K136").Value = 3429.15
$ws.Range("L136").Value = 70203.429
$ws.Range("M136").Value = -879.1499999999996
$ws.Range("N136").Value = -75303.429

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2502.2222
$ws.Range("I105").Value = 2505
$ws.Range("J105").Value = 2480
$ws.Range("K105").Value = 2505
$ws.Range("L105").Value = 2480
$ws.Range("M105").Value = -758
$ws.Range("N105").Value = -5974
# Row 134
$ws.Range("H134").Value = 773290.1
$ws.Range("I134").Value = 1084874
$ws.Range("J134").Value = 4716.467
$ws.Range("K134").Value = 3254622
$ws.Range("L134").Value = 14149.401
$ws.Range("M134").Value = -3252087
$ws.Range("N134").Value = -19219.401
# Row 135
$ws.Range("H135").Value = 25470.588
$ws.Range("J135").Value = 25470.588
$ws.Range("L135").Value = 25470.588
$ws.Range("N135").Value = -35610.588

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2350.4167
$ws.Range("J62").Value = 2350
$ws.Range("L62").Value = 2350
$ws.Range("N62").Value = -3598
# Row 65
$ws.Range("H65").Value = 2350.4167
$ws.Range("J65").Value = 2350
$ws.Range("L65").Value = 11750
$ws.Range("N65").Value = -17990
# Row 82
$ws.Range("H82").Value = 28110.924
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 28786.834
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 28786.834
$ws.Range("M82").Value = -19639
$ws.Range("N82").Value = -29508.834
# Row 85
$ws.Range("H85").Value = 28110.924
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 28786.834
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 28786.834
$ws.Range("M85").Value = -18752
$ws.Range("N85").Value = -31282.834
# Row 99
$ws.Range("H99").Value = 6466.6665
$ws.Range("I99").Value = 6700
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 6700
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -5202
$ws.Range("N99").Value = -8996
# Row 122
$ws.Range("H122").Value = 937.3333
$ws.Range("I122").Value = 937.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2811.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -361.9998999999998
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 6466.6665
$ws.Range("I126").Value = 6700
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 20100
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -17630
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2289.7368
$ws.Range("I80").Value = 2275.4167
$ws.Range("J80").Value = 2314.2856
$ws.Range("K80").Value = 2275.4167
$ws.Range("L80").Value = 2314.2856
$ws.Range("M80").Value = -1277.4167
$ws.Range("N80").Value = -4310.2856
# Row 83
$ws.Range("H83").Value = 2289.7368
$ws.Range("I83").Value = 2275.4167
$ws.Range("J83").Value = 2314.2856
$ws.Range("K83").Value = 11377.0835
$ws.Range("L83").Value = 11571.428
$ws.Range("M83").Value = -6385.083500000001
$ws.Range("N83").Value = -21555.428

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# Row 40
$ws.Range("H40").Value = 7166.6665
$ws.Range("I40").Value = 7166.6665
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7166.6665
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7030.6665
$ws.Range("N40").ClearContents()
# Row 68
$ws.Range("H68").Value = 26459.375
$ws.Range("I68").Value = 67167.336
$ws.Range("J68").Value = 2034.6
$ws.Range("K68").Value = 67167.336
$ws.Range("L68").Value = 2034.6
$ws.Range("M68").Value = -66418.336
$ws.Range("N68").Value = -3532.6
# Row 71
$ws.Range("H71").Value = 26459.375
$ws.Range("I71").Value = 67167.336
$ws.Range("J71").Value = 2034.6
$ws.Range("K71").Value = 335836.68
$ws.Range("L71").Value = 10173
$ws.Range("M71").Value = -332092.68
$ws.Range("N71").Value = -17661
# Row 82
$ws.Range("H82").Value = 891.6539
$ws.Range("I82").Value = 814
$ws.Range("J82").Value = 948.6
$ws.Range("K82").Value = 814
$ws.Range("L82").Value = 948.6
$ws.Range("M82").Value = -453
$ws.Range("N82").Value = -1670.6
# Row 85
$ws.Range("H85").Value = 891.6539
$ws.Range("I85").Value = 814
$ws.Range("J85").Value = 948.6
$ws.Range("K85").Value = 814
$ws.Range("L85").Value = 948.6
$ws.Range("M85").Value = 434
$ws.Range("N85").Value = -3444.6
# Row 100
$ws.Range("H100").Value = 1526.4736
$ws.Range("I100").Value = 1384.8462
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 1384.8462
$ws.Range("L100").Value = 1833.3334
$ws.Range("M100").Value = -843.8462
$ws.Range("N100").Value = -2915.3334
# Row 122
$ws.Range("H122").Value = 3444.2856
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3444.2856
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10332.8568
$ws.Range("N122").Value = -15232.8568
$ws.Range("M122").ClearContents()
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 134501.11
$ws.Range("I122").Value = 833.3333
$ws.Range("K122").Value = 2499.9999
$ws.Range("M122").Value = -49.9998999999998
